$d = $word.ActiveDocument

function Set-ParaInnerXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $r.MoveEnd(1, -1)
    $r.Text = ""
    $head = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>'
    $tail = '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $xml = $head + $innerXml + $tail
    $r.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the end of "Basic memory hacking app."
#    to the new TODO item about std::remove_cv (added below).
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Insert a new list item after "Merge headers where appropriate ..."
#    containing the note about std::remove_cv, with the _GoBack bookmark
#    placed between the "std::" run and the final run.
# ---------------------------------------------------------------------------
$mergeParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -match "^Merge headers where appropriate") {
        $mergeParaIndex = $i
        break
    }
}

$mergePara = $d.Paragraphs.Item($mergeParaIndex)
$mergeRange = $mergePara.Range
$mergeRange.Collapse(0)
$mergeRange.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($mergeParaIndex + 1)
$newRange = $newPara.Range
$newRange.MoveEnd(1, -1)

$newInner = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Check whether Read/Write APIs and other templates should be using </w:t></w:r><w:r><w:t>std::</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>remove_cv anywhere in the type detection/transformation.</w:t></w:r>'
$newHead = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>'
$newTail = '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newXml = $newHead + $newInner + $newTail
$newRange.InsertXML($newXml)

Write-Output "Inserted new TODO item after paragraph $mergeParaIndex."

# ---------------------------------------------------------------------------
# 3) Re-paginate the "lastRenderedPageBreak" markers: each one moves to the
#    run that now starts the printed page after the new paragraph shifted
#    everything down.
# ---------------------------------------------------------------------------
function Find-ParaIndex($pattern) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs.Item($i).Range.Text
        if ($t -match $pattern) {
            return $i
        }
    }
    return -1
}

$memoryMgrIdx = Find-ParaIndex "^MemoryMgr"
$uncheckedIdx = Find-ParaIndex "^.Unchecked. read/write"
$detectIdx = Find-ParaIndex "^Detect cases where hooking"
$findPatternIdx = Find-ParaIndex "^FindPattern"
$nopIdx = Find-ParaIndex "^NOP/UnNOP support"
$scannerIdx = Find-ParaIndex "^Scanner"

Write-Output "memoryMgrIdx=$memoryMgrIdx uncheckedIdx=$uncheckedIdx detectIdx=$detectIdx findPatternIdx=$findPatternIdx nopIdx=$nopIdx scannerIdx=$scannerIdx"

# MemoryMgr: add lastRenderedPageBreak before the bold run text.
$memMgrInner = '<w:pPr><w:ind w:left="2160" w:hanging="2160"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:lastRenderedPageBreak/><w:t>MemoryMgr</w:t></w:r>'
Set-ParaInnerXml $memoryMgrIdx $memMgrInner

# 'Unchecked' read/write ...: remove lastRenderedPageBreak.
$uncheckedInner = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>' + [char]0x2018 + 'Unchecked' + [char]0x2019 + ' read/write etc functions designed for speed and use in ReadString etc where you only want to check page protections once, then forget about it.</w:t></w:r>'
Set-ParaInnerXml $uncheckedIdx $uncheckedInner

# Detect cases where hooking ...: add lastRenderedPageBreak.
$detectInner = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Detect cases where hooking may overflow past the end of a function, and fail. (Provide policy or flag to allow overriding this behaviour.) Examples may be instructions such as ' + [char]0x2018 + 'int 3' + [char]0x2019 + ', ' + [char]0x2018 + 'ret' + [char]0x2019 + ', ' + [char]0x2018 + 'jmp' + [char]0x2019 + ', etc.</w:t></w:r>'
Set-ParaInnerXml $detectIdx $detectInner

# FindPattern: remove lastRenderedPageBreak.
$findPatternInner = '<w:r><w:rPr><w:b/></w:rPr><w:t>FindPattern</w:t></w:r>'
Set-ParaInnerXml $findPatternIdx $findPatternInner

# NOP/UnNOP support.: add lastRenderedPageBreak.
$nopInner = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>NOP/UnNOP support.</w:t></w:r>'
Set-ParaInnerXml $nopIdx $nopInner

# Scanner: remove lastRenderedPageBreak.
$scannerInner = '<w:r><w:rPr><w:b/></w:rPr><w:t>Scanner</w:t></w:r>'
Set-ParaInnerXml $scannerIdx $scannerInner

Write-Output "Done."
